$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.371.92"
$ws.Range("E2").Value = "  +5.74%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.557.77"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "412.95"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.80"
$ws.Range("E6").Value = "  -1.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.643"
$ws.Range("E7").Value = "  +3.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.546.02"
$ws.Range("E8").Value = "  +2.58%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.767"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.172"
$ws.Range("E11").Value = "  +18.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000312"
$ws.Range("E12").Value = "  +43.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.93"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.77"
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.113.58"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.687.70"
$ws.Range("E17").Value = "  +6.95%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.07"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.352.93"
$ws.Range("E20").Value = "  +6.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.20"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "443.83"
$ws.Range("E22").Value = "  -4.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.69"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.12"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.92"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.29"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.93"
$ws.Range("E27").Value = "  -6.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.34"
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.76"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.23"
$ws.Range("E31").Value = "  +2.03%  "
$ws.Range("E32").Value = "  +3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.25"
$ws.Range("E33").Value = "  -4.61%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.32"
$ws.Range("E36").Value = "  -3.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.31"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0486"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0717"
$ws.Range("E39").Value = "  +29.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.146"
$ws.Range("E40").Value = "  +9.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.96"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "146.78"
$ws.Range("E44").Value = "  +0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.22"
$ws.Range("E45").Value = "  -3.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.28"
$ws.Range("E46").Value = "  -1.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.305"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.95"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.25"
$ws.Range("E49").Value = "  -6.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.00"
$ws.Range("E50").Value = "  +6.16%  "
$ws.Range("B51").Value = "Celestia"
$ws.Range("C51").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "15.40"
$ws.Range("E51").Value = "  -5.75%  "
